$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two data rows were dropped from the "missing data" sample entirely
# (not just blanked): "RM 232" (old row 26) and "SC 92" (old row 28,
# which is row 27 once RM 232 is already gone). Deleting them shifts
# everything below up accordingly.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# The remaining rows also got a different "missingness" pattern applied
# to column B (the "A" measurement):
#   - "SC 101" (now row 27) previously had no value there -> now filled in.
#   - "SC 119" (now row 29) previously had a value -> now blanked out.
#   - "SC 193" (now row 32) previously had a value -> now blanked out.
$ws.Cells.Item(27, 2).Value = -20.4
$ws.Cells.Item(29, 2).ClearContents()
$ws.Cells.Item(32, 2).ClearContents()
